$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values that changed
$ws.Cells.Item(21, 2).Value = 146321
$ws.Cells.Item(22, 2).Value = 164191
$ws.Cells.Item(23, 2).Value = 177668
$ws.Cells.Item(29, 2).Value = 14032
$ws.Cells.Item(30, 2).Value = 16000
$ws.Cells.Item(31, 2).Value = 15655

# Append new Senate election rows (37-60)
$newRows = @(
    @{ Row = 37; Label = "Senate Election 1996"; Value = 588 }
    @{ Row = 38; Label = "Senate Election 1998"; Value = 140 }
    @{ Row = 39; Label = "Senate By-election 1999"; Value = 8 }
    @{ Row = 40; Label = "Senate Election 2000"; Value = 165 }
    @{ Row = 41; Label = "Senate Election 2002"; Value = 199 }
    @{ Row = 42; Label = "Senate By-election 2003"; Value = 16 }
    @{ Row = 43; Label = "Senate By-election 2004"; Value = 12 }
    @{ Row = 44; Label = "Senate Election 2004"; Value = 216 }
    @{ Row = 45; Label = "Senate Election 2006"; Value = 214 }
    @{ Row = 46; Label = "Senate By-election 2007"; Value = 17 }
    @{ Row = 47; Label = "Senate Election 2008"; Value = 203 }
    @{ Row = 48; Label = "Senate Election 2010"; Value = 228 }
    @{ Row = 49; Label = "Senate By-election 2011"; Value = 10 }
    @{ Row = 50; Label = "Senate Election 2012"; Value = 236 }
    @{ Row = 51; Label = "Senate By-election 2014"; Value = 20 }
    @{ Row = 52; Label = "Senate Election 2014"; Value = 245 }
    @{ Row = 53; Label = "Senate By-election 2016"; Value = 234 }
    @{ Row = 54; Label = "Senate By-election 2017"; Value = 9 }
    @{ Row = 55; Label = "Senate By-election 2018"; Value = 18 }
    @{ Row = 56; Label = "Senate Election 2018"; Value = 237 }
    @{ Row = 57; Label = "Senate By-election 2019"; Value = 10 }
    @{ Row = 58; Label = "Senate By-election 2020"; Value = 11 }
    @{ Row = 59; Label = "Senate Election 2020"; Value = 237 }
    @{ Row = 60; Label = "Senate Election 2022"; Value = 181 }
)

foreach ($item in $newRows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Label
    $ws.Cells.Item($item.Row, 2).Value = $item.Value
}
